# "Updated documentation to show remaining tasks"
#
# Mark a set of backlog tasks as "Complete" in column D (the Status
# column) now that they have been finished. Each row number below
# corresponds to a task listed in column C of the "Backlog" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$completedRows = @(12, 32, 34, 35, 38, 39, 40, 72, 75, 78, 79, 80)

foreach ($row in $completedRows) {
    $ws.Cells.Item($row, 4).Value = "Complete"
}
